$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Elimina EC anteriores y se agregan nuevos" -- the 7 previous "Periodo Mora"
# account-statement periods (1905..1911) are replaced by the new periods in
# reverse order (1911..1905), and the "Valor Mora" amounts for the first and
# last worker rows are swapped to match the re-ordered periods.
$ws.Range("E16").Value = "1911"
$ws.Range("E17").Value = "1910"
$ws.Range("E18").Value = "1909"
$ws.Range("E19").Value = "1908"
$ws.Range("E20").Value = "1907"
$ws.Range("E21").Value = "1906"
$ws.Range("E22").Value = "1905"

$ws.Range("F16").Value = 6944
$ws.Range("F22").Value = 3168

# Se modifica base de datos: columnas ligeramente mas anchas para ajustarse
# al nuevo contenido (autofit recompute under the new Excel build).
$ws.Columns.Item(2).ColumnWidth = 17.7
$ws.Columns.Item(3).ColumnWidth = 15.9
$ws.Columns.Item(5).ColumnWidth = 12.71
$ws.Columns.Item(6).ColumnWidth = 9.35
$ws.Columns.Item(7).ColumnWidth = 13.5
$ws.Columns.Item(8).ColumnWidth = 18.5
$ws.Columns.Item(9).ColumnWidth = 17.33
$ws.Columns.Item(10).ColumnWidth = 14.17
